$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (31) following the same pattern as row 30, describing
# the new "simulacion" feature (SimulacionController / SP_DIVIDENDOSHIS_SELECT).
$ws.Range("A31").Value = "simulacion"
$ws.Range("B31").Value = "SimulacionController"
$ws.Range("C31").Value = "SP_DIVIDENDOSHIS_SELECT"
$ws.Range("D31").Value = '$emisor,$precio,$capita'
$ws.Range("G31").Value = "shares-page"

# Copy styles from row 30 so formatting matches the rest of the table.
$ws.Range("A30").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B30").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null

$ws.Range("C30").Copy() | Out-Null
$ws.Range("C31").PasteSpecial(-4122) | Out-Null

$ws.Range("G30").Copy() | Out-Null
$ws.Range("G31").PasteSpecial(-4122) | Out-Null

# Update the active selection to match the post-edit state recorded in the diff.
$ws.Range("C37").Select() | Out-Null
